# Applies a shuffle of the species-observation records in rows 15-38
# (row 24 is unchanged) on the "Artfynd" worksheet. Each destination row
# ends up holding the values (for the columns that actually differ) that
# used to live in a different source row; the row numbers themselves
# never move, only the data originally found in them.
#
# destRow -> srcRow mapping used below:
# 15<-20 16<-25 17<-29 18<-22 19<-35 20<-18 21<-33 22<-30 23<-28 25<-15 26<-23 27<-36 28<-27 29<-17 30<-38 31<-37 32<-19 33<-16 34<-21 35<-26 36<-31 37<-34 38<-32

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15
$ws.Cells.Item(15,1).Value2 = 111667773
$ws.Cells.Item(15,2).Value2 = 89790
$ws.Cells.Item(15,4).Value2 = "NT"
$ws.Cells.Item(15,5).Value2 = 6040186
$ws.Cells.Item(15,6).ClearContents() | Out-Null
$ws.Cells.Item(15,7).Value2 = "Leptoporus mollis"
$ws.Cells.Item(15,8).Value2 = "(Pers.:Fr.) Quél."
$ws.Cells.Item(15,10).ClearContents() | Out-Null
$ws.Cells.Item(15,12).ClearContents() | Out-Null
$ws.Cells.Item(15,13).ClearContents() | Out-Null
$ws.Cells.Item(15,14).ClearContents() | Out-Null
$ws.Cells.Item(15,16).Value2 = "Svartflärken (Svartflärken), Mpd"
$ws.Cells.Item(15,17).Value2 = 628028.5695976926
$ws.Cells.Item(15,18).Value2 = 6944391.799111729
$ws.Cells.Item(15,32).ClearContents() | Out-Null

# Row 16
$ws.Cells.Item(16,1).Value2 = 111667338
$ws.Cells.Item(16,2).Value2 = 89425
$ws.Cells.Item(16,4).Value2 = "NT"
$ws.Cells.Item(16,5).Value2 = 5442
$ws.Cells.Item(16,6).Value2 = "Tallticka"
$ws.Cells.Item(16,7).Value2 = "Porodaedalea pini"
$ws.Cells.Item(16,8).Value2 = "(Brot.) Murrill"
$ws.Cells.Item(16,10).ClearContents() | Out-Null
$ws.Cells.Item(16,12).ClearContents() | Out-Null
$ws.Cells.Item(16,13).ClearContents() | Out-Null
$ws.Cells.Item(16,14).ClearContents() | Out-Null
$ws.Cells.Item(16,16).Value2 = "Bladbacken (Bladbacken), Mpd"
$ws.Cells.Item(16,17).Value2 = 628206.6965146795
$ws.Cells.Item(16,18).Value2 = 6944047.015808343
$ws.Cells.Item(16,32).ClearContents() | Out-Null

# Row 17
$ws.Cells.Item(17,1).Value2 = 111667041
$ws.Cells.Item(17,2).Value2 = 77550
$ws.Cells.Item(17,4).Value2 = "NT"
$ws.Cells.Item(17,5).Value2 = 185
$ws.Cells.Item(17,6).Value2 = "Violettgrå tagellav"
$ws.Cells.Item(17,7).Value2 = "Bryoria nadvornikiana"
$ws.Cells.Item(17,8).Value2 = "(Gyeln.) Brodo & D.Hawksw."
$ws.Cells.Item(17,10).ClearContents() | Out-Null
$ws.Cells.Item(17,12).ClearContents() | Out-Null
$ws.Cells.Item(17,13).ClearContents() | Out-Null
$ws.Cells.Item(17,14).ClearContents() | Out-Null
$ws.Cells.Item(17,16).Value2 = "Svartflärksbäcken (Svartflärksbäcken), Mpd"
$ws.Cells.Item(17,17).Value2 = 627994.5336875709
$ws.Cells.Item(17,18).Value2 = 6943827.688791481
$ws.Cells.Item(17,32).ClearContents() | Out-Null

# Row 18
$ws.Cells.Item(18,1).Value2 = 111667763
$ws.Cells.Item(18,2).Value2 = 89845
$ws.Cells.Item(18,4).Value2 = "VU"
$ws.Cells.Item(18,5).Value2 = 1209
$ws.Cells.Item(18,6).Value2 = "Rynkskinn"
$ws.Cells.Item(18,7).Value2 = "Phlebia centrifuga"
$ws.Cells.Item(18,8).Value2 = "P.Karst."
$ws.Cells.Item(18,10).ClearContents() | Out-Null
$ws.Cells.Item(18,12).ClearContents() | Out-Null
$ws.Cells.Item(18,13).ClearContents() | Out-Null
$ws.Cells.Item(18,14).ClearContents() | Out-Null
$ws.Cells.Item(18,16).Value2 = "Svartflärken (Svartflärken), Mpd"
$ws.Cells.Item(18,17).Value2 = 628028.5695976926
$ws.Cells.Item(18,18).Value2 = 6944391.799111729
$ws.Cells.Item(18,32).ClearContents() | Out-Null

# Row 19
$ws.Cells.Item(19,1).Value2 = 111666918
$ws.Cells.Item(19,2).Value2 = 77550
$ws.Cells.Item(19,4).Value2 = "NT"
$ws.Cells.Item(19,5).Value2 = 185
$ws.Cells.Item(19,6).Value2 = "Violettgrå tagellav"
$ws.Cells.Item(19,7).Value2 = "Bryoria nadvornikiana"
$ws.Cells.Item(19,8).Value2 = "(Gyeln.) Brodo & D.Hawksw."
$ws.Cells.Item(19,10).ClearContents() | Out-Null
$ws.Cells.Item(19,12).ClearContents() | Out-Null
$ws.Cells.Item(19,13).ClearContents() | Out-Null
$ws.Cells.Item(19,14).ClearContents() | Out-Null
$ws.Cells.Item(19,16).Value2 = "Svartflärksbäcken (Svartflärksbäcken), Mpd"
$ws.Cells.Item(19,17).Value2 = 627981.5521892406
$ws.Cells.Item(19,18).Value2 = 6943733.609182604
$ws.Cells.Item(19,32).ClearContents() | Out-Null

# Row 20
$ws.Cells.Item(20,1).Value2 = 111667287
$ws.Cells.Item(20,2).Value2 = 78578
$ws.Cells.Item(20,4).Value2 = "NT"
$ws.Cells.Item(20,5).Value2 = 6458
$ws.Cells.Item(20,6).Value2 = "Lunglav"
$ws.Cells.Item(20,7).Value2 = "Lobaria pulmonaria"
$ws.Cells.Item(20,8).Value2 = "(L.) Hoffm."
$ws.Cells.Item(20,10).ClearContents() | Out-Null
$ws.Cells.Item(20,12).ClearContents() | Out-Null
$ws.Cells.Item(20,13).ClearContents() | Out-Null
$ws.Cells.Item(20,14).ClearContents() | Out-Null
$ws.Cells.Item(20,16).Value2 = "Bladbacken (Bladbacken), Mpd"
$ws.Cells.Item(20,17).Value2 = 628191.351562822
$ws.Cells.Item(20,18).Value2 = 6944014.155575473
$ws.Cells.Item(20,32).ClearContents() | Out-Null

# Row 21
$ws.Cells.Item(21,1).Value2 = 111668151
$ws.Cells.Item(21,2).Value2 = 78578
$ws.Cells.Item(21,4).Value2 = "NT"
$ws.Cells.Item(21,5).Value2 = 6458
$ws.Cells.Item(21,6).Value2 = "Lunglav"
$ws.Cells.Item(21,7).Value2 = "Lobaria pulmonaria"
$ws.Cells.Item(21,8).Value2 = "(L.) Hoffm."
$ws.Cells.Item(21,10).ClearContents() | Out-Null
$ws.Cells.Item(21,12).ClearContents() | Out-Null
$ws.Cells.Item(21,13).ClearContents() | Out-Null
$ws.Cells.Item(21,14).ClearContents() | Out-Null
$ws.Cells.Item(21,16).Value2 = "Svartflärken (Svartflärken), Mpd"
$ws.Cells.Item(21,17).Value2 = 627992.8558976713
$ws.Cells.Item(21,18).Value2 = 6944372.443055111
$ws.Cells.Item(21,32).ClearContents() | Out-Null

# Row 22
$ws.Cells.Item(22,1).Value2 = 111668301
$ws.Cells.Item(22,2).Value2 = 77550
$ws.Cells.Item(22,4).Value2 = "NT"
$ws.Cells.Item(22,5).Value2 = 185
$ws.Cells.Item(22,6).Value2 = "Violettgrå tagellav"
$ws.Cells.Item(22,7).Value2 = "Bryoria nadvornikiana"
$ws.Cells.Item(22,8).Value2 = "(Gyeln.) Brodo & D.Hawksw."
$ws.Cells.Item(22,10).ClearContents() | Out-Null
$ws.Cells.Item(22,12).ClearContents() | Out-Null
$ws.Cells.Item(22,13).ClearContents() | Out-Null
$ws.Cells.Item(22,14).ClearContents() | Out-Null
$ws.Cells.Item(22,16).Value2 = "Svartflärken (Svartflärken), Mpd"
$ws.Cells.Item(22,17).Value2 = 627869.5672010599
$ws.Cells.Item(22,18).Value2 = 6944134.919311633
$ws.Cells.Item(22,32).ClearContents() | Out-Null

# Row 23
$ws.Cells.Item(23,1).Value2 = 111668569
$ws.Cells.Item(23,2).Value2 = 77550
$ws.Cells.Item(23,4).Value2 = "NT"
$ws.Cells.Item(23,5).Value2 = 185
$ws.Cells.Item(23,6).Value2 = "Violettgrå tagellav"
$ws.Cells.Item(23,7).Value2 = "Bryoria nadvornikiana"
$ws.Cells.Item(23,8).Value2 = "(Gyeln.) Brodo & D.Hawksw."
$ws.Cells.Item(23,10).ClearContents() | Out-Null
$ws.Cells.Item(23,12).ClearContents() | Out-Null
$ws.Cells.Item(23,13).ClearContents() | Out-Null
$ws.Cells.Item(23,14).ClearContents() | Out-Null
$ws.Cells.Item(23,16).Value2 = "Svartflärksbäcken (Svartflärksbäcken), Mpd"
$ws.Cells.Item(23,17).Value2 = 627931.2258218131
$ws.Cells.Item(23,18).Value2 = 6943686.026669092
$ws.Cells.Item(23,32).ClearContents() | Out-Null

# Row 25
$ws.Cells.Item(25,1).Value2 = 111668313
$ws.Cells.Item(25,2).Value2 = 56543
$ws.Cells.Item(25,4).Value2 = "NT"
$ws.Cells.Item(25,5).Value2 = 103021
$ws.Cells.Item(25,6).Value2 = "Talltita"
$ws.Cells.Item(25,7).Value2 = "Poecile montanus"
$ws.Cells.Item(25,8).Value2 = "(Conrad von Baldenstein, 1827)"
$ws.Cells.Item(25,10).ClearContents() | Out-Null
$ws.Cells.Item(25,12).ClearContents() | Out-Null
$ws.Cells.Item(25,13).ClearContents() | Out-Null
$ws.Cells.Item(25,14).ClearContents() | Out-Null
$ws.Cells.Item(25,16).Value2 = "Svartflärken (Svartflärken), Mpd"
$ws.Cells.Item(25,17).Value2 = 627869.5672010599
$ws.Cells.Item(25,18).Value2 = 6944134.919311633
$ws.Cells.Item(25,32).ClearContents() | Out-Null

# Row 26
$ws.Cells.Item(26,1).Value2 = 111667281
$ws.Cells.Item(26,2).Value2 = 85715
$ws.Cells.Item(26,4).Value2 = "NT"
$ws.Cells.Item(26,5).Value2 = 510
$ws.Cells.Item(26,6).Value2 = "Doftskinn"
$ws.Cells.Item(26,7).Value2 = "Cystostereum murrayi"
$ws.Cells.Item(26,8).Value2 = "(Berk. & M.A. Curtis.) Pouzar"
$ws.Cells.Item(26,10).ClearContents() | Out-Null
$ws.Cells.Item(26,12).ClearContents() | Out-Null
$ws.Cells.Item(26,13).ClearContents() | Out-Null
$ws.Cells.Item(26,14).ClearContents() | Out-Null
$ws.Cells.Item(26,16).Value2 = "Bladbacken (Bladbacken), Mpd"
$ws.Cells.Item(26,17).Value2 = 628196.0116143352
$ws.Cells.Item(26,18).Value2 = 6944024.937745438
$ws.Cells.Item(26,32).ClearContents() | Out-Null

# Row 27
$ws.Cells.Item(27,1).Value2 = 111667927
$ws.Cells.Item(27,2).Value2 = 77550
$ws.Cells.Item(27,4).Value2 = "NT"
$ws.Cells.Item(27,5).Value2 = 185
$ws.Cells.Item(27,6).Value2 = "Violettgrå tagellav"
$ws.Cells.Item(27,7).Value2 = "Bryoria nadvornikiana"
$ws.Cells.Item(27,8).Value2 = "(Gyeln.) Brodo & D.Hawksw."
$ws.Cells.Item(27,10).ClearContents() | Out-Null
$ws.Cells.Item(27,12).ClearContents() | Out-Null
$ws.Cells.Item(27,13).ClearContents() | Out-Null
$ws.Cells.Item(27,14).ClearContents() | Out-Null
$ws.Cells.Item(27,16).Value2 = "Svartflärken (Svartflärken), Mpd"
$ws.Cells.Item(27,17).Value2 = 628070.1241137966
$ws.Cells.Item(27,18).Value2 = 6944749.195084839
$ws.Cells.Item(27,32).ClearContents() | Out-Null

# Row 28
$ws.Cells.Item(28,1).Value2 = 111667731
$ws.Cells.Item(28,2).Value2 = 77515
$ws.Cells.Item(28,4).Value2 = "NT"
$ws.Cells.Item(28,5).Value2 = 6425
$ws.Cells.Item(28,6).Value2 = "Garnlav"
$ws.Cells.Item(28,7).Value2 = "Alectoria sarmentosa"
$ws.Cells.Item(28,8).Value2 = "(Ach.) Ach."
$ws.Cells.Item(28,10).ClearContents() | Out-Null
$ws.Cells.Item(28,12).ClearContents() | Out-Null
$ws.Cells.Item(28,13).ClearContents() | Out-Null
$ws.Cells.Item(28,14).ClearContents() | Out-Null
$ws.Cells.Item(28,16).Value2 = "Svartflärken (Svartflärken), Mpd"
$ws.Cells.Item(28,17).Value2 = 628005.9220808987
$ws.Cells.Item(28,18).Value2 = 6944356.817101943
$ws.Cells.Item(28,32).ClearContents() | Out-Null

# Row 29
$ws.Cells.Item(29,1).Value2 = 111667725
$ws.Cells.Item(29,2).Value2 = 77550
$ws.Cells.Item(29,4).Value2 = "NT"
$ws.Cells.Item(29,5).Value2 = 185
$ws.Cells.Item(29,6).Value2 = "Violettgrå tagellav"
$ws.Cells.Item(29,7).Value2 = "Bryoria nadvornikiana"
$ws.Cells.Item(29,8).Value2 = "(Gyeln.) Brodo & D.Hawksw."
$ws.Cells.Item(29,10).ClearContents() | Out-Null
$ws.Cells.Item(29,12).ClearContents() | Out-Null
$ws.Cells.Item(29,13).ClearContents() | Out-Null
$ws.Cells.Item(29,14).ClearContents() | Out-Null
$ws.Cells.Item(29,16).Value2 = "Svartflärken (Svartflärken), Mpd"
$ws.Cells.Item(29,17).Value2 = 628005.9220808987
$ws.Cells.Item(29,18).Value2 = 6944356.817101943
$ws.Cells.Item(29,32).ClearContents() | Out-Null

# Row 30
$ws.Cells.Item(30,1).Value2 = 111668169
$ws.Cells.Item(30,2).Value2 = 89405
$ws.Cells.Item(30,4).Value2 = "NT"
$ws.Cells.Item(30,5).Value2 = 1202
$ws.Cells.Item(30,6).Value2 = "Ullticka"
$ws.Cells.Item(30,7).Value2 = "Phellinidium ferrugineofuscum"
$ws.Cells.Item(30,8).Value2 = "(P.Karst.) Fiasson & Niemelä"
$ws.Cells.Item(30,10).ClearContents() | Out-Null
$ws.Cells.Item(30,12).ClearContents() | Out-Null
$ws.Cells.Item(30,13).ClearContents() | Out-Null
$ws.Cells.Item(30,14).ClearContents() | Out-Null
$ws.Cells.Item(30,16).Value2 = "Svartflärken (Svartflärken), Mpd"
$ws.Cells.Item(30,17).Value2 = 627989.210128115
$ws.Cells.Item(30,18).Value2 = 6944335.430490699
$ws.Cells.Item(30,32).ClearContents() | Out-Null

# Row 31
$ws.Cells.Item(31,1).Value2 = 111667427
$ws.Cells.Item(31,2).Value2 = 77550
$ws.Cells.Item(31,4).Value2 = "NT"
$ws.Cells.Item(31,5).Value2 = 185
$ws.Cells.Item(31,6).Value2 = "Violettgrå tagellav"
$ws.Cells.Item(31,7).Value2 = "Bryoria nadvornikiana"
$ws.Cells.Item(31,8).Value2 = "(Gyeln.) Brodo & D.Hawksw."
$ws.Cells.Item(31,10).ClearContents() | Out-Null
$ws.Cells.Item(31,12).ClearContents() | Out-Null
$ws.Cells.Item(31,13).ClearContents() | Out-Null
$ws.Cells.Item(31,14).ClearContents() | Out-Null
$ws.Cells.Item(31,16).Value2 = "Bladbacken (Bladbacken), Mpd"
$ws.Cells.Item(31,17).Value2 = 628238.2251299906
$ws.Cells.Item(31,18).Value2 = 6944162.548277185
$ws.Cells.Item(31,32).ClearContents() | Out-Null

# Row 32
$ws.Cells.Item(32,1).Value2 = 111667877
$ws.Cells.Item(32,2).Value2 = 56543
$ws.Cells.Item(32,4).Value2 = "NT"
$ws.Cells.Item(32,5).Value2 = 103021
$ws.Cells.Item(32,6).Value2 = "Talltita"
$ws.Cells.Item(32,7).Value2 = "Poecile montanus"
$ws.Cells.Item(32,8).Value2 = "(Conrad von Baldenstein, 1827)"
$ws.Cells.Item(32,10).ClearContents() | Out-Null
$ws.Cells.Item(32,12).ClearContents() | Out-Null
$ws.Cells.Item(32,13).ClearContents() | Out-Null
$ws.Cells.Item(32,14).ClearContents() | Out-Null
$ws.Cells.Item(32,16).Value2 = "Svartflärken (Svartflärken), Mpd"
$ws.Cells.Item(32,17).Value2 = 628073.1002953692
$ws.Cells.Item(32,18).Value2 = 6944660.366756786
$ws.Cells.Item(32,32).ClearContents() | Out-Null

# Row 33
$ws.Cells.Item(33,1).Value2 = 111667778
$ws.Cells.Item(33,2).Value2 = 77550
$ws.Cells.Item(33,4).Value2 = "NT"
$ws.Cells.Item(33,5).Value2 = 185
$ws.Cells.Item(33,6).Value2 = "Violettgrå tagellav"
$ws.Cells.Item(33,7).Value2 = "Bryoria nadvornikiana"
$ws.Cells.Item(33,8).Value2 = "(Gyeln.) Brodo & D.Hawksw."
$ws.Cells.Item(33,10).ClearContents() | Out-Null
$ws.Cells.Item(33,12).ClearContents() | Out-Null
$ws.Cells.Item(33,13).ClearContents() | Out-Null
$ws.Cells.Item(33,14).ClearContents() | Out-Null
$ws.Cells.Item(33,16).Value2 = "Svartflärken (Svartflärken), Mpd"
$ws.Cells.Item(33,17).Value2 = 628041.295580395
$ws.Cells.Item(33,18).Value2 = 6944396.900308819
$ws.Cells.Item(33,32).ClearContents() | Out-Null

# Row 34
$ws.Cells.Item(34,1).Value2 = 111667231
$ws.Cells.Item(34,2).Value2 = 77550
$ws.Cells.Item(34,4).Value2 = "NT"
$ws.Cells.Item(34,5).Value2 = 185
$ws.Cells.Item(34,6).Value2 = "Violettgrå tagellav"
$ws.Cells.Item(34,7).Value2 = "Bryoria nadvornikiana"
$ws.Cells.Item(34,8).Value2 = "(Gyeln.) Brodo & D.Hawksw."
$ws.Cells.Item(34,10).ClearContents() | Out-Null
$ws.Cells.Item(34,12).ClearContents() | Out-Null
$ws.Cells.Item(34,13).ClearContents() | Out-Null
$ws.Cells.Item(34,14).ClearContents() | Out-Null
$ws.Cells.Item(34,16).Value2 = "Bladbacken (Bladbacken), Mpd"
$ws.Cells.Item(34,17).Value2 = 628160.9420623753
$ws.Cells.Item(34,18).Value2 = 6944048.470567195
$ws.Cells.Item(34,32).ClearContents() | Out-Null

# Row 35
$ws.Cells.Item(35,1).Value2 = 111667471
$ws.Cells.Item(35,2).Value2 = 77515
$ws.Cells.Item(35,4).Value2 = "NT"
$ws.Cells.Item(35,5).Value2 = 6425
$ws.Cells.Item(35,6).Value2 = "Garnlav"
$ws.Cells.Item(35,7).Value2 = "Alectoria sarmentosa"
$ws.Cells.Item(35,8).Value2 = "(Ach.) Ach."
$ws.Cells.Item(35,10).ClearContents() | Out-Null
$ws.Cells.Item(35,12).ClearContents() | Out-Null
$ws.Cells.Item(35,13).ClearContents() | Out-Null
$ws.Cells.Item(35,14).ClearContents() | Out-Null
$ws.Cells.Item(35,16).Value2 = "Svartflärken (Svartflärken), Mpd"
$ws.Cells.Item(35,17).Value2 = 628242.7901880945
$ws.Cells.Item(35,18).Value2 = 6944199.598880037
$ws.Cells.Item(35,32).ClearContents() | Out-Null

# Row 36
$ws.Cells.Item(36,1).Value2 = 111667081
$ws.Cells.Item(36,2).Value2 = 77550
$ws.Cells.Item(36,4).Value2 = "NT"
$ws.Cells.Item(36,5).Value2 = 185
$ws.Cells.Item(36,6).Value2 = "Violettgrå tagellav"
$ws.Cells.Item(36,7).Value2 = "Bryoria nadvornikiana"
$ws.Cells.Item(36,8).Value2 = "(Gyeln.) Brodo & D.Hawksw."
$ws.Cells.Item(36,10).ClearContents() | Out-Null
$ws.Cells.Item(36,12).ClearContents() | Out-Null
$ws.Cells.Item(36,13).ClearContents() | Out-Null
$ws.Cells.Item(36,14).ClearContents() | Out-Null
$ws.Cells.Item(36,16).Value2 = "Svartflärksbäcken (Svartflärksbäcken), Mpd"
$ws.Cells.Item(36,17).Value2 = 628030.2196984198
$ws.Cells.Item(36,18).Value2 = 6943835.983260213
$ws.Cells.Item(36,32).ClearContents() | Out-Null

# Row 37
$ws.Cells.Item(37,1).Value2 = 111668109
$ws.Cells.Item(37,2).Value2 = 5135
$ws.Cells.Item(37,4).Value2 = "LC"
$ws.Cells.Item(37,5).Value2 = 105930
$ws.Cells.Item(37,6).Value2 = "Vågbandad barkbock"
$ws.Cells.Item(37,7).Value2 = "Semanotus undatus"
$ws.Cells.Item(37,8).Value2 = "(Linnaeus, 1758)"
$ws.Cells.Item(37,10).ClearContents() | Out-Null
$ws.Cells.Item(37,12).ClearContents() | Out-Null
$ws.Cells.Item(37,13).Value2 = "äldre gnagspår"
$ws.Cells.Item(37,14).ClearContents() | Out-Null
$ws.Cells.Item(37,16).Value2 = "Svartflärken (Svartflärken), Mpd"
$ws.Cells.Item(37,17).Value2 = 628016.3600143436
$ws.Cells.Item(37,18).Value2 = 6944480.738271755
$ws.Cells.Item(37,32).ClearContents() | Out-Null

# Row 38
$ws.Cells.Item(38,1).Value2 = 111667977
$ws.Cells.Item(38,2).Value2 = 77550
$ws.Cells.Item(38,4).Value2 = "NT"
$ws.Cells.Item(38,5).Value2 = 185
$ws.Cells.Item(38,6).Value2 = "Violettgrå tagellav"
$ws.Cells.Item(38,7).Value2 = "Bryoria nadvornikiana"
$ws.Cells.Item(38,8).Value2 = "(Gyeln.) Brodo & D.Hawksw."
$ws.Cells.Item(38,10).ClearContents() | Out-Null
$ws.Cells.Item(38,12).ClearContents() | Out-Null
$ws.Cells.Item(38,13).ClearContents() | Out-Null
$ws.Cells.Item(38,14).ClearContents() | Out-Null
$ws.Cells.Item(38,16).Value2 = "Svartflärken (Svartflärken), Mpd"
$ws.Cells.Item(38,17).Value2 = 628025.0361159063
$ws.Cells.Item(38,18).Value2 = 6944745.60709906
$ws.Cells.Item(38,32).ClearContents() | Out-Null

